$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.656.60'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.845.02'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4323'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3706'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07342'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8789'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.99'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '1.852.82'
$ws.Range("E12").Value = '  -2.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.483'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.607'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009039'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '27.559.21'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.134'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.74%  '
$ws.Range("D24").Value = '2.144.46'
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.992'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.875'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08938'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7875'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.619'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.177'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.976'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05445'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.104'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01965'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.853'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5181'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1691'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.793'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.641'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4793'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06556'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.666'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.841'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.18%  '
